$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 256.8
$ws.Range("I12").Value = 164.5
$ws.Range("J12").Value = 395.25
$ws.Range("K12").Value = 164.5
$ws.Range("L12").Value = 395.25
$ws.Range("M12").Value = 5.5
$ws.Range("N12").Value = -735.25
$ws.Range("H19").Value = 941.7857
$ws.Range("J19").Value = 911.1111
$ws.Range("L19").Value = 911.1111
$ws.Range("N19").Value = -1261.1111
$ws.Range("H28").Value = 4946.375
$ws.Range("I28").Value = 447.27274
$ws.Range("K28").Value = 447.27274
$ws.Range("M28").Value = 37.72726
$ws.Range("H62").Value = 71429950
$ws.Range("I62").Value = 83334670
$ws.Range("K62").Value = 83334670
$ws.Range("M62").Value = -83334046
$ws.Range("H65").Value = 71429950
$ws.Range("I65").Value = 83334670
$ws.Range("K65").Value = 416673350
$ws.Range("M65").Value = -416670230
$ws.Range("H70").Value = 5438.8887
$ws.Range("I70").Value = 9000
$ws.Range("J70").Value = 4993.75
$ws.Range("K70").Value = 27000
$ws.Range("L70").Value = 14981.25
$ws.Range("M70").Value = -26730
$ws.Range("N70").Value = -15521.25
$ws.Range("H73").Value = 5438.8887
$ws.Range("I73").Value = 9000
$ws.Range("J73").Value = 4993.75
$ws.Range("K73").Value = 27000
$ws.Range("L73").Value = 14981.25
$ws.Range("M73").Value = -26064
$ws.Range("N73").Value = -16853.25
$ws.Range("H107").Value = 35717284
$ws.Range("I107").Value = 62500720
$ws.Range("K107").Value = 62500720
$ws.Range("M107").Value = -62498800
$ws.Range("H112").Value = 2599.6
$ws.Range("J112").Value = 2682.3794
$ws.Range("L112").Value = 8047.138199999999
$ws.Range("N112").Value = -10263.1382
$ws.Range("H136").Value = 67250
$ws.Range("J136").Value = 84500
$ws.Range("L136").Value = 84500
$ws.Range("N136").Value = -94700
$ws.Range("H137").Value = 4391510.5
$ws.Range("I137").Value = 6946955
$ws.Range("K137").Value = 20840865
$ws.Range("M137").Value = -20838315
$ws.Range("H138").Value = 4080.4385
$ws.Range("I138").Value = 2554.3125
$ws.Range("J138").Value = 4676
$ws.Range("K138").Value = 7662.9375
$ws.Range("L138").Value = 14028
$ws.Range("M138").Value = -2522.9375
$ws.Range("N138").Value = -24308

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4907173.5
$ws.Range("I32").Value = 6095070
$ws.Range("J32").Value = 7099.25
$ws.Range("K32").Value = 6095070
$ws.Range("L32").Value = 7099.25
$ws.Range("M32").Value = -6094783
$ws.Range("N32").Value = -7673.25
$ws.Range("H61").Value = 6935.533
$ws.Range("I61").Value = 4852
$ws.Range("K61").Value = 4852
$ws.Range("M61").Value = -4640
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""
$ws.Range("H136").Value = 6935.533
$ws.Range("I136").Value = 4852
$ws.Range("K136").Value = 14556
$ws.Range("M136").Value = -12006

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2702.5557
$ws.Range("I107").Value = 1656.3549
$ws.Range("J107").Value = 9189
$ws.Range("K107").Value = 1656.3549
$ws.Range("L107").Value = 9189
$ws.Range("M107").Value = 263.6451
$ws.Range("N107").Value = -13029
$ws.Range("H134").Value = 6485.7354
$ws.Range("I134").Value = 2802.875
$ws.Range("J134").Value = 9759.388999999999
$ws.Range("K134").Value = 8408.625
$ws.Range("L134").Value = 29278.167
$ws.Range("M134").Value = -5873.625
$ws.Range("N134").Value = -34348.167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4467.3
$ws.Range("I31").Value = 2276
$ws.Range("K31").Value = 2276
$ws.Range("M31").Value = -1981
$ws.Range("H34").Value = 4467.3
$ws.Range("I34").Value = 2276
$ws.Range("K34").Value = 2276
$ws.Range("M34").Value = -2074
$ws.Range("H58").Value = 4943.7407
$ws.Range("I58").Value = 3109.2778
$ws.Range("K58").Value = 3109.2778
$ws.Range("M58").Value = -2906.2778
$ws.Range("H86").Value = 7500
$ws.Range("I86").Value = 5832
$ws.Range("J86").Value = 12504
$ws.Range("K86").Value = 5832
$ws.Range("L86").Value = 12504
$ws.Range("M86").Value = -4709
$ws.Range("N86").Value = -14750
$ws.Range("H89").Value = 7500
$ws.Range("I89").Value = 5832
$ws.Range("J89").Value = 12504
$ws.Range("K89").Value = 29160
$ws.Range("L89").Value = 62520
$ws.Range("M89").Value = -23544
$ws.Range("N89").Value = -73752
$ws.Range("H109").Value = 79990
$ws.Range("J109").Value = 79990
$ws.Range("L109").Value = 79990
$ws.Range("N109").Value = -82070
$ws.Range("H134").Value = 5496.727
$ws.Range("I134").Value = 5226.4287
$ws.Range("K134").Value = 15679.2861
$ws.Range("M134").Value = -13144.2861
$ws.Range("H136").Value = 4943.7407
$ws.Range("I136").Value = 3109.2778
$ws.Range("K136").Value = 9327.8334
$ws.Range("M136").Value = -6777.8334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 10164.818
$ws.Range("J11").Value = 56
$ws.Range("L11").Value = 168
$ws.Range("N11").Value = -448
$ws.Range("H68").Value = 46003.652
$ws.Range("I68").Value = 201189.9
$ws.Range("J68").Value = 2896.361
$ws.Range("K68").Value = 603569.7
$ws.Range("L68").Value = 8689.082999999999
$ws.Range("M68").Value = -602758.7
$ws.Range("N68").Value = -10311.083
$ws.Range("H71").Value = 46003.652
$ws.Range("I71").Value = 201189.9
$ws.Range("J71").Value = 2896.361
$ws.Range("K71").Value = 1810709.1
$ws.Range("L71").Value = 26067.249
$ws.Range("M71").Value = -1806653.1
$ws.Range("N71").Value = -34179.249
$ws.Range("H141").Value = 1890
$ws.Range("I141").Value = 1890
$ws.Range("K141").Value = 5670
$ws.Range("M141").Value = -490

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 10052.333
$ws.Range("I132").Value = 4362.4
$ws.Range("K132").Value = 13087.2
$ws.Range("M132").Value = -10557.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5227.9443
$ws.Range("I7").Value = 5938.857
$ws.Range("J7").Value = 2739.75
$ws.Range("K7").Value = 5938.857
$ws.Range("L7").Value = 2739.75
$ws.Range("M7").Value = -5826.857
$ws.Range("N7").Value = -2963.75
$ws.Range("H22").Value = 2353.2068
$ws.Range("I22").Value = 1805.8966
$ws.Range("J22").Value = 2900.5173
$ws.Range("K22").Value = 1805.8966
$ws.Range("L22").Value = 2900.5173
$ws.Range("M22").Value = -1510.8966
$ws.Range("N22").Value = -3490.5173
$ws.Range("H27").Value = 2353.2068
$ws.Range("I27").Value = 1805.8966
$ws.Range("J27").Value = 2900.5173
$ws.Range("K27").Value = 1805.8966
$ws.Range("L27").Value = 2900.5173
$ws.Range("M27").Value = -1698.8966
$ws.Range("N27").Value = -3114.5173
$ws.Range("H55").Value = 599.8095
$ws.Range("I55").Value = 335.5
$ws.Range("J55").Value = 1128.4286
$ws.Range("K55").Value = 335.5
$ws.Range("L55").Value = 1128.4286
$ws.Range("M55").Value = -162.5
$ws.Range("N55").Value = -1474.4286
$ws.Range("H61").Value = 4906.9414
$ws.Range("I61").Value = 4906.9414
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4906.9414
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4704.9414
$ws.Range("N61").Value = ""
$ws.Range("H113").Value = 4906.9414
$ws.Range("I113").Value = 4906.9414
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4906.9414
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2736.9414
$ws.Range("N113").Value = ""
$ws.Range("H126").Value = 5227.9443
$ws.Range("I126").Value = 5938.857
$ws.Range("J126").Value = 2739.75
$ws.Range("K126").Value = 17816.571
$ws.Range("L126").Value = 8219.25
$ws.Range("M126").Value = -15346.571
$ws.Range("N126").Value = -13159.25
$ws.Range("H132").Value = 6648.5713
$ws.Range("I132").Value = 4308.1
$ws.Range("J132").Value = 12499.75
$ws.Range("K132").Value = 12924.3
$ws.Range("L132").Value = 37499.25
$ws.Range("M132").Value = -10394.3
$ws.Range("N132").Value = -42559.25
$ws.Range("H136").Value = 5469.0386
$ws.Range("I136").Value = 3324.9375
$ws.Range("K136").Value = 9974.8125
$ws.Range("M136").Value = -7424.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 5359.6113
$ws.Range("I61").Value = 5492.933
$ws.Range("J61").Value = 4693
$ws.Range("K61").Value = 5492.933
$ws.Range("L61").Value = 4693
$ws.Range("M61").Value = -5200.933
$ws.Range("N61").Value = -5277
$ws.Range("H107").Value = 3308.879
$ws.Range("I107").Value = 3415.2222
$ws.Range("J107").Value = 2830.3333
$ws.Range("K107").Value = 10245.6666
$ws.Range("L107").Value = 8490.999899999999
$ws.Range("M107").Value = -8325.6666
$ws.Range("N107").Value = -12330.9999
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = ""
$ws.Range("H113").Value = 443.33334
$ws.Range("I113").Value = 461.66666
$ws.Range("K113").Value = 1384.99998
$ws.Range("M113").Value = 785.0000199999999
$ws.Range("H132").Value = 3873.7192
$ws.Range("I132").Value = 2135.535
$ws.Range("K132").Value = 6406.605
$ws.Range("M132").Value = -3876.605
